# Fix never being able to validate
# - Strip leading spaces from header names and rename a few columns
# - Add a new "ContactName" column (U) populated with "plk"
# - Normalize diseaseTested values to "Bd"
# - Normalize sampleMethod/sampleType values to "external Swab"
# - Fix fieldNumber sequence for the last three rows (7, 8, 9)
# - Switch page orientation to portrait
# - Leave selection on D9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header renames (also strips the leading space the old sheet had) ---
$ws.Range("B1").Value = "catalogNumber"
$ws.Range("C1").Value = "fieldNumber"
$ws.Range("D1").Value = "diseaseTested"
$ws.Range("E1").Value = "diseaseStrain"
$ws.Range("F1").Value = "sampleType"
$ws.Range("G1").Value = "sampleDisposition"
$ws.Range("H1").Value = "diseaseDetected"
$ws.Range("I1").Value = "fatal"
$ws.Range("J1").Value = "cladeSampled"
$ws.Range("K1").Value = "genus"
$ws.Range("L1").Value = "specificEpithet"
$ws.Range("M1").Value = "infraspecificEpithet"
$ws.Range("N1").Value = "lifeStage"
$ws.Range("O1").Value = "dateCollected"
$ws.Range("P1").Value = "decimalLatitude"
$ws.Range("Q1").Value = "decimalLongitude"
$ws.Range("R1").Value = "elevation"
$ws.Range("S1").Value = "coordinateUncertaintyInMeters"
$ws.Range("T1").Value = "Collector"
$ws.Range("U1").Value = "ContactName"

# --- Data rows 2-10 ---
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "Bd"              # D: diseaseTested
    $ws.Cells.Item($r, 6).Value = "external Swab"   # F: sampleType
    $ws.Cells.Item($r, 21).Value = "plk"             # U: ContactName
}

# fieldNumber (column C) sequence fix for the last three rows
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 8
$ws.Range("C10").Value = 9

# --- Page setup / selection ---
$ws.PageSetup.Orientation = 1
$ws.Range("D9").Select() | Out-Null
